# Adds a new "2/15/20" data column (AC) to the right of the existing
# "2/14/20" column (AB) on the COVID19-Deaths sheet, mirroring the source
# CSV-to-Excel export appending one more day of data.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# --- Header cell AC1: "2/15/20" -------------------------------------------
# A plain Value assignment of a date-like string ("2/15/20") gets
# auto-recognized by Excel as a date serial number. Prefixing with a
# leading apostrophe forces it to be stored as literal text (matching the
# other date-header cells, which are plain shared strings).
$ws.Range("AC1").Value = "'2/15/20"

# Match the bold/centered/bordered header formatting used by the rest of
# row 1 (same style as the neighboring AB1 cell).
$ws.Range("AB1").Copy()
$ws.Range("AC1").PasteSpecial(-4122)

# --- Data cells AC2:AC76: death counts for 2/15/20 ------------------------
$deaths = @(6,4,5,0,2,2,2,1,4,3,11,13,1596,2,0,0,1,1,1,0,0,0,2,1,0,1,3,0,1,0,0,0,1,0,0,0,0,0,0,1,0,0,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,1,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0,0)

$rows = $deaths.Length
$arr = New-Object 'object[,]' $rows,1
for ($i = 0; $i -lt $rows; $i++) {
    $arr[$i,0] = $deaths[$i]
}

$ws.Range("AC2:AC76").Value = $arr
